$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.179.76'
$ws.Range('E2').Value = '  -3.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.807.46'
$ws.Range('E3').Value = '  -3.79%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.86'
$ws.Range('E5').Value = '  -2.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4192'
$ws.Range('E7').Value = '  -2.72%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3553'
$ws.Range('E8').Value = '  -4.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07110'
$ws.Range('E9').Value = '  -4.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8469'
$ws.Range('E10').Value = '  -4.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.15'
$ws.Range('E11').Value = '  -4.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.771.40'
$ws.Range('E12').Value = '  -7.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.294'
$ws.Range('E13').Value = '  -3.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.354'
$ws.Range('E14').Value = '  -4.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06854'
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.005'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '80.69'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008748'
$ws.Range('E18').Value = '  -4.08%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.12'
$ws.Range('E20').Value = '  -3.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.323.21'
$ws.Range('E21').Value = '  -3.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.093'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.81'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.046.31'
$ws.Range('E24').Value = '  -8.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.967'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.59'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.13'
$ws.Range('E27').Value = '  -3.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.067'
$ws.Range('E28').Value = '  -6.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.86'
$ws.Range('E29').Value = '  -3.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.675'
$ws.Range('E30').Value = '  -10.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08894'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.949'
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7325'
$ws.Range('E33').Value = '  -7.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.374'
$ws.Range('E34').Value = '  -6.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.101'
$ws.Range('E35').Value = '  -6.75%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.079'
$ws.Range('E37').Value = '  -4.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05147'
$ws.Range('E38').Value = '  -5.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01900'
$ws.Range('E39').Value = '  -3.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.703'
$ws.Range('E40').Value = '  -6.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1626'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4957'
$ws.Range('E42').Value = '  -4.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.240'
$ws.Range('E43').Value = '  -9.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.138'
$ws.Range('E44').Value = '  -6.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.68'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.23'
$ws.Range('E46').Value = '  -3.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06351'
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4555'
$ws.Range('E49').Value = '  -4.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.594'
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.25'
$ws.Range('E51').Value = '  -5.26%  '
